$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new "Zugnummer" column before current column C (Datum shifts C->D) ---
$ws.Columns("C:C").Insert()
$ws.Range("C5").Value = "Zugnummer"

# --- Split "Gastfahrt vor Dienstbeginn" (now at E5) into "... von" / "... bis" ---
$ws.Range("E5").Value = "Gastfahrt vor Dienstbeginn von"
$ws.Columns("F:F").Insert()
$ws.Range("F5").Value = "Gastfahrt vor Dienstbeginn bis"

# --- Split "Pause" (now at J5) into "Pause von" / "Pause bis" ---
$ws.Range("J5").Value = "Pause von"
$ws.Columns("K:K").Insert()
$ws.Range("K5").Value = "Pause bis"

# --- Split "Wartezeit 1" (now at L5) into "Wartezeit 1 von" / "Wartezeit 1 bis" ---
$ws.Range("L5").Value = "Wartezeit 1 von"
$ws.Columns("M:M").Insert()
$ws.Range("M5").Value = "Wartezeit 1 bis"

# --- Split "Wartezeit 2" (now at N5) into "Wartezeit 2 von" / "Wartezeit 2 bis" ---
$ws.Range("N5").Value = "Wartezeit 2 von"
$ws.Columns("O:O").Insert()
$ws.Range("O5").Value = "Wartezeit 2 bis"

# --- Split "Abfahrt / Ankunft" (now at S5) into "Abfahrt" / "Ankunft" ---
$ws.Range("S5").Value = "Abfahrt"
$ws.Columns("T:T").Insert()
$ws.Range("T5").Value = "Ankunft"

# --- Split "Gastfahrt nach Dienstende" (now at V5) into "... von" / "... bis" ---
$ws.Range("V5").Value = "Gastfahrt nach Dienstende von"
$ws.Columns("W:W").Insert()
$ws.Range("W5").Value = "Gastfahrt nach Dienstende bis"

# --- Restore the active cell selection to C5 ---
$ws.Range("C5").Select() | Out-Null
